# Mandatory RCS and Kerbal GPS Revived Support
# Update MonoPropellant engine upgrade inputs (Thrust / ISP source values)
# and refresh the active selection, matching the authored workbook edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the base input values that drive the rest of the sheet's formulas ---
$ws.Range("B2").Value = 560    # Thrust
$ws.Range("C2").Value = 337    # ISP
$ws.Range("C3").Value = 213    # ISP (Methalox row)
$ws.Range("B7").Value = 13000  # Upgrade Cost

# --- Move the active selection to where the author left it ---
$ws.Range("D3").Select()
